$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 697, shifting the current
# rows 697:720 down to 700:723 (dimension grows from T720 to T723).
$ws.Range("A697:A699").EntireRow.Insert()

# --- New row 697: Platano "Pinton", week of 2022-07-?? (serial 44747) ---
$ws.Cells.Item(697, 1).Value = 8
$ws.Cells.Item(697, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(697, 3).Value = "Coquimbo"
$ws.Cells.Item(697, 4).Value = 44747
$ws.Cells.Item(697, 5).Value = 4
$ws.Cells.Item(697, 6).Value = "Fruta"
$ws.Cells.Item(697, 7).Value = 100108
$ws.Cells.Item(697, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(697, 9).Value = 100108006
$ws.Cells.Item(697, 10).Value = "Plátano"
$ws.Cells.Item(697, 11).Value = "Sin especificar"
$ws.Cells.Item(697, 12).Value = "Pintón"
$ws.Cells.Item(697, 13).Value = 120
$ws.Cells.Item(697, 14).Value = 23000
$ws.Cells.Item(697, 15).Value = 23000
$ws.Cells.Item(697, 16).Value = 23000
$ws.Cells.Item(697, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(697, 18).Value = "Ecuador"
$ws.Cells.Item(697, 19).Value = 1150
$ws.Cells.Item(697, 20).Value = 20

# --- New row 698: Platano "Primera Maduro" ---
$ws.Cells.Item(698, 1).Value = 8
$ws.Cells.Item(698, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(698, 3).Value = "Coquimbo"
$ws.Cells.Item(698, 4).Value = 44747
$ws.Cells.Item(698, 5).Value = 4
$ws.Cells.Item(698, 6).Value = "Fruta"
$ws.Cells.Item(698, 7).Value = 100108
$ws.Cells.Item(698, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(698, 9).Value = 100108006
$ws.Cells.Item(698, 10).Value = "Plátano"
$ws.Cells.Item(698, 11).Value = "Sin especificar"
$ws.Cells.Item(698, 12).Value = "Primera Maduro"
$ws.Cells.Item(698, 13).Value = 120
$ws.Cells.Item(698, 14).Value = 25000
$ws.Cells.Item(698, 15).Value = 25000
$ws.Cells.Item(698, 16).Value = 25000
$ws.Cells.Item(698, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(698, 18).Value = "Ecuador"
$ws.Cells.Item(698, 19).Value = 1250
$ws.Cells.Item(698, 20).Value = 20

# --- New row 699: Platano "Primera Pinton" ---
$ws.Cells.Item(699, 1).Value = 8
$ws.Cells.Item(699, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(699, 3).Value = "Coquimbo"
$ws.Cells.Item(699, 4).Value = 44747
$ws.Cells.Item(699, 5).Value = 4
$ws.Cells.Item(699, 6).Value = "Fruta"
$ws.Cells.Item(699, 7).Value = 100108
$ws.Cells.Item(699, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(699, 9).Value = 100108006
$ws.Cells.Item(699, 10).Value = "Plátano"
$ws.Cells.Item(699, 11).Value = "Sin especificar"
$ws.Cells.Item(699, 12).Value = "Primera Pintón"
$ws.Cells.Item(699, 13).Value = 160
$ws.Cells.Item(699, 14).Value = 26000
$ws.Cells.Item(699, 15).Value = 26000
$ws.Cells.Item(699, 16).Value = 26000
$ws.Cells.Item(699, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(699, 18).Value = "Ecuador"
$ws.Cells.Item(699, 19).Value = 1300
$ws.Cells.Item(699, 20).Value = 20
